$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (introduces shared strings 0-6 in column order)
$ws.Range("A1").Value = "petId"
$ws.Range("B1").Value = "categoryId"
$ws.Range("C1").Value = "categoryName"
$ws.Range("D1").Value = "petName"
$ws.Range("E1").Value = "tagId"
$ws.Range("F1").Value = "tagName"
$ws.Range("G1").Value = "status"

# Row 2 - first pet (Bruno). Note petName (D2) is written before
# categoryName (C2) so the shared-string table matches the source order
# (index 7 = Bruno, index 8 = Labrador).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = "Bruno"
$ws.Range("C2").Value = "Labrador"
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = "bruno"
$ws.Range("G2").Value = "available"

# Row 3 - second pet (Milo)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = "German Shephard"
$ws.Range("D3").Value = "Milo"
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = "milo"
$ws.Range("G3").Value = "available"

# Match the final cursor/selection shown in the diff.
$ws.Range("C11").Select()

# Auto-fit the columns so the stored widths reflect the new content,
# matching the bestFit column definitions in the target sheet.
$ws.Columns("A:G").AutoFit()
